## Generate Report for Handoff
## Adds two new localization entries (6b7ed46d-... and debb9713-...) to the
## Overview, zh-cn and de-de sheets, mirroring the existing "Ready for
## handoff" rows (row 3 in each sheet).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview": columns A..G, header in row1, data rows 2-3 existing.
# Append rows 4 and 5 (copy format from row 3, then set new values).
# ---------------------------------------------------------------------
$wsOv = $wb.Worksheets.Item("Overview")

$wsOv.Rows(3).Copy()
$wsOv.Rows(4).Insert(-4121)
$wsOv.Rows(3).Copy()
$wsOv.Rows(5).Insert(-4121)

$wsOv.Range("A4").Value = "6b7ed46d-656e-4c8c-9ed1-7cf5aa61d787.md"
$wsOv.Range("B4").Value = "e2e\6b7ed46d-656e-4c8c-9ed1-7cf5aa61d787.md"
$wsOv.Range("C4").Value = ".md"
$wsOv.Range("D4").Value = ""
$wsOv.Range("E4").Value = "Ready for handoff"
$wsOv.Range("F4").Value = "Ready for handoff"
$wsOv.Range("G4").Value = "2016-08-13 12:51:17"

$wsOv.Range("A5").Value = "debb9713-69cf-4251-82b8-e29b77304eea.md"
$wsOv.Range("B5").Value = "e2e\debb9713-69cf-4251-82b8-e29b77304eea.md"
$wsOv.Range("C5").Value = ".md"
$wsOv.Range("D5").Value = ""
$wsOv.Range("E5").Value = "Ready for handoff"
$wsOv.Range("F5").Value = "Ready for handoff"
$wsOv.Range("G5").Value = "2016-08-13 12:51:17"

$wsOv.Hyperlinks.Add($wsOv.Range("B4"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/304c33bb9dfb9a77ba49ba749b8d94cc57179f8e/e2e/6b7ed46d-656e-4c8c-9ed1-7cf5aa61d787.md", $null, $null, "e2e\6b7ed46d-656e-4c8c-9ed1-7cf5aa61d787.md") | Out-Null
$wsOv.Hyperlinks.Add($wsOv.Range("B5"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/304c33bb9dfb9a77ba49ba749b8d94cc57179f8e/e2e/debb9713-69cf-4251-82b8-e29b77304eea.md", $null, $null, "e2e\debb9713-69cf-4251-82b8-e29b77304eea.md") | Out-Null

$loOv = $wsOv.ListObjects.Item(1)
$loOv.Resize($wsOv.Range("A1:G5"))

# ---------------------------------------------------------------------
# Sheet "zh-cn": columns A..P, header row1, data rows 2-3 existing.
# Row 3 is the template ("Ready for handoff", no I/J hyperlink).
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Rows(3).Copy()
$wsZh.Rows(4).Insert(-4121)
$wsZh.Rows(3).Copy()
$wsZh.Rows(5).Insert(-4121)

$wsZh.Range("A4").Value = "6b7ed46d-656e-4c8c-9ed1-7cf5aa61d787.md"
$wsZh.Range("B4").Value = ".md"
$wsZh.Range("C4").Value = "Ready for handoff"
$wsZh.Range("D4").Value = "e2e"
$wsZh.Range("E4").Value = "ht"
$wsZh.Range("F4").Value = "'False"
$wsZh.Range("G4").Value = "6b7ed46d-656e-4c8c-9ed1-7cf5aa61d787.8149e0e5135ebe5ceb6371e10bd546667f8d6c64.zh-cn.xlf"
$wsZh.Range("H4").Value = "2016-08-13 12:51:10"
$wsZh.Range("I4").Value = ""
$wsZh.Range("J4").Value = ""
$wsZh.Range("K4").Value = "0001-01-01 00:00:00"
$wsZh.Range("L4").Value = ""
$wsZh.Range("M4").Value = "'True"
$wsZh.Range("N4").Value = ""
$wsZh.Range("O4").Value = "'False"
$wsZh.Range("P4").Value = ""

$wsZh.Range("A5").Value = "debb9713-69cf-4251-82b8-e29b77304eea.md"
$wsZh.Range("B5").Value = ".md"
$wsZh.Range("C5").Value = "Ready for handoff"
$wsZh.Range("D5").Value = "e2e"
$wsZh.Range("E5").Value = "ht"
$wsZh.Range("F5").Value = "'False"
$wsZh.Range("G5").Value = "debb9713-69cf-4251-82b8-e29b77304eea.13b06a3efcac90dfa8e0dfeb15c996bdd4ecd8dc.zh-cn.xlf"
$wsZh.Range("H5").Value = "2016-08-13 12:51:10"
$wsZh.Range("I5").Value = ""
$wsZh.Range("J5").Value = ""
$wsZh.Range("K5").Value = "0001-01-01 00:00:00"
$wsZh.Range("L5").Value = ""
$wsZh.Range("M5").Value = "'True"
$wsZh.Range("N5").Value = ""
$wsZh.Range("O5").Value = "'False"
$wsZh.Range("P5").Value = ""

$wsZh.Hyperlinks.Add($wsZh.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/95047f6bbea44f79790eff81c4878e1cf5309cc1/e2e/6b7ed46d-656e-4c8c-9ed1-7cf5aa61d787.md", $null, $null, "6b7ed46d-656e-4c8c-9ed1-7cf5aa61d787.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("A5"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/95047f6bbea44f79790eff81c4878e1cf5309cc1/e2e/debb9713-69cf-4251-82b8-e29b77304eea.md", $null, $null, "debb9713-69cf-4251-82b8-e29b77304eea.md") | Out-Null

$loZh = $wsZh.ListObjects.Item(1)
$loZh.Resize($wsZh.Range("A1:P5"))

# ---------------------------------------------------------------------
# Sheet "de-de": columns A..P, header row1, data rows 2-3 existing.
# Row 3 is the template ("Ready for handoff", no I/J hyperlink).
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Rows(3).Copy()
$wsDe.Rows(4).Insert(-4121)
$wsDe.Rows(3).Copy()
$wsDe.Rows(5).Insert(-4121)

$wsDe.Range("A4").Value = "6b7ed46d-656e-4c8c-9ed1-7cf5aa61d787.md"
$wsDe.Range("B4").Value = ".md"
$wsDe.Range("C4").Value = "Ready for handoff"
$wsDe.Range("D4").Value = "e2e"
$wsDe.Range("E4").Value = "ht"
$wsDe.Range("F4").Value = "'False"
$wsDe.Range("G4").Value = "6b7ed46d-656e-4c8c-9ed1-7cf5aa61d787.8149e0e5135ebe5ceb6371e10bd546667f8d6c64.de-de.xlf"
$wsDe.Range("H4").Value = "2016-08-13 12:51:17"
$wsDe.Range("I4").Value = ""
$wsDe.Range("J4").Value = ""
$wsDe.Range("K4").Value = "0001-01-01 00:00:00"
$wsDe.Range("L4").Value = ""
$wsDe.Range("M4").Value = "'True"
$wsDe.Range("N4").Value = ""
$wsDe.Range("O4").Value = "'False"
$wsDe.Range("P4").Value = ""

$wsDe.Range("A5").Value = "debb9713-69cf-4251-82b8-e29b77304eea.md"
$wsDe.Range("B5").Value = ".md"
$wsDe.Range("C5").Value = "Ready for handoff"
$wsDe.Range("D5").Value = "e2e"
$wsDe.Range("E5").Value = "ht"
$wsDe.Range("F5").Value = "'False"
$wsDe.Range("G5").Value = "debb9713-69cf-4251-82b8-e29b77304eea.13b06a3efcac90dfa8e0dfeb15c996bdd4ecd8dc.de-de.xlf"
$wsDe.Range("H5").Value = "2016-08-13 12:51:17"
$wsDe.Range("I5").Value = ""
$wsDe.Range("J5").Value = ""
$wsDe.Range("K5").Value = "0001-01-01 00:00:00"
$wsDe.Range("L5").Value = ""
$wsDe.Range("M5").Value = "'True"
$wsDe.Range("N5").Value = ""
$wsDe.Range("O5").Value = "'False"
$wsDe.Range("P5").Value = ""

$wsDe.Hyperlinks.Add($wsDe.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/9bf7ca5f519e2848eb08401f70598e920b08ae02/e2e/6b7ed46d-656e-4c8c-9ed1-7cf5aa61d787.md", $null, $null, "6b7ed46d-656e-4c8c-9ed1-7cf5aa61d787.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("A5"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/9bf7ca5f519e2848eb08401f70598e920b08ae02/e2e/debb9713-69cf-4251-82b8-e29b77304eea.md", $null, $null, "debb9713-69cf-4251-82b8-e29b77304eea.md") | Out-Null

$loDe = $wsDe.ListObjects.Item(1)
$loDe.Resize($wsDe.Range("A1:P5"))
